$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15, shifting the Category Default block (old rows 15-28) down to 16-29
$ws.Rows("15:15").Insert()

# Populate the new row 15 with the isInitialization config fields
$ws.Range("C15").Value = "isInitialization"
$ws.Range("D15").Value = "BOOL"
$ws.Range("E15").Value = "0- not initialization /1 - initialization"

# Match the thin/medium top border treatment used by the other "header divider" rows
$ws.Range("C15:E15").Borders.Item(8).Weight = 2

$ws.Range("C15").Borders.Item(7).Weight = -4138
$ws.Range("C15").Borders.Item(10).Weight = 2

$ws.Range("D15").Borders.Item(7).Weight = 2
$ws.Range("D15").Borders.Item(10).Weight = 2

$ws.Range("E15").Borders.Item(7).Weight = 2
$ws.Range("E15").Borders.Item(10).Weight = -4138

# Note explaining how isInitialization interacts with CategoryID defaulting
$ws.Range("E11").Value = "If initialization is 1 => catID = 1"

# Widen column E so the longer note text fits
$ws.Range("E1").ColumnWidth = 30.5546875

Write-Output "edit applied"
